# Auto-generated COM script applying scheduled price-data refresh to Aegis_Profits workbook
# (unique diff: 44 row updates across 8 sheets; values are plain numeric market-price caches, no formulas)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 756.2
$ws.Range("I107").Value = 718.05884
$ws.Range("J107").Value = 837.25
$ws.Range("K107").Value = 718.05884
$ws.Range("L107").Value = 837.25
$ws.Range("M107").Value = 1201.94116
$ws.Range("N107").Value = -4677.25

$ws.Range("H112").Value = 1723.1666
$ws.Range("J112").Value = 2007.8
$ws.Range("L112").Value = 6023.4
$ws.Range("N112").Value = -8239.4

$ws.Range("H123").Value = 17962
$ws.Range("J123").Value = 17962
$ws.Range("L123").Value = 17962
$ws.Range("N123").Value = -27762

$ws.Range("H129").Value = 7544.8
$ws.Range("I129").Value = 17192.5
$ws.Range("J129").Value = 1113
$ws.Range("K129").Value = 51577.5
$ws.Range("L129").Value = 3339
$ws.Range("M129").Value = -46577.5
$ws.Range("N129").Value = -13339

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null

$ws.Range("H137").Value = 1624.5834
$ws.Range("I137").Value = 1633.6923
$ws.Range("J137").Value = 1613.8182
$ws.Range("K137").Value = 4901.0769
$ws.Range("L137").Value = 4841.4546
$ws.Range("M137").Value = -2351.0769
$ws.Range("N137").Value = -9941.454600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 4777.778
$ws.Range("I3").Value = 3000
$ws.Range("K3").Value = 3000
$ws.Range("M3").Value = -2885

$ws.Range("H32").Value = 38731.547
$ws.Range("I32").Value = 17330.867
$ws.Range("J32").Value = 124334.266
$ws.Range("K32").Value = 17330.867
$ws.Range("L32").Value = 124334.266
$ws.Range("M32").Value = -17043.867
$ws.Range("N32").Value = -124908.266

$ws.Range("H63").Value = 3100
$ws.Range("J63").Value = 3100
$ws.Range("L63").Value = 3100
$ws.Range("N63").Value = -4472

$ws.Range("H66").Value = 3100
$ws.Range("J66").Value = 3100
$ws.Range("L66").Value = 15500
$ws.Range("N66").Value = -22364

$ws.Range("H74").Value = 1643.6522
$ws.Range("I74").Value = 1579
$ws.Range("J74").Value = 1744.2222
$ws.Range("K74").Value = 1579
$ws.Range("L74").Value = 1744.2222
$ws.Range("M74").Value = -705
$ws.Range("N74").Value = -3492.2222

$ws.Range("H77").Value = 1643.6522
$ws.Range("I77").Value = 1579
$ws.Range("J77").Value = 1744.2222
$ws.Range("K77").Value = 7895
$ws.Range("L77").Value = 8721.110999999999
$ws.Range("M77").Value = -3527
$ws.Range("N77").Value = -17457.111

$ws.Range("H97").Value = 48897.383
$ws.Range("I97").Value = 53674.895
$ws.Range("K97").Value = 53674.895
$ws.Range("M97").Value = -53178.895

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 3005
$ws.Range("I12").Value = 3005
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 3005
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -2837
$ws.Range("N12").Value = $null

$ws.Range("H99").Value = 1664.2858
$ws.Range("I99").Value = 1386.1538
$ws.Range("J99").Value = 2116.25
$ws.Range("K99").Value = 1386.1538
$ws.Range("L99").Value = 2116.25
$ws.Range("M99").Value = 111.8462
$ws.Range("N99").Value = -5112.25

$ws.Range("H127").Value = 38822.5
$ws.Range("J127").Value = 38822.5
$ws.Range("L127").Value = 38822.5
$ws.Range("N127").Value = -48742.5

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 681.5
$ws.Range("I22").Value = 195.66667
$ws.Range("J22").Value = 1167.3334
$ws.Range("K22").Value = 195.66667
$ws.Range("L22").Value = 1167.3334
$ws.Range("M22").Value = 154.33333
$ws.Range("N22").Value = -1867.3334

$ws.Range("H31").Value = 23259.03
$ws.Range("I31").Value = 1108
$ws.Range("J31").Value = 74944.766
$ws.Range("K31").Value = 1108
$ws.Range("L31").Value = 74944.766
$ws.Range("M31").Value = -813
$ws.Range("N31").Value = -75534.766

$ws.Range("H34").Value = 23259.03
$ws.Range("I34").Value = 1108
$ws.Range("J34").Value = 74944.766
$ws.Range("K34").Value = 1108
$ws.Range("L34").Value = 74944.766
$ws.Range("M34").Value = -906
$ws.Range("N34").Value = -75348.766

$ws.Range("H50").Value = 9457.6
$ws.Range("J50").Value = 9457.6
$ws.Range("L50").Value = 9457.6
$ws.Range("N50").Value = -10707.6

$ws.Range("H51").Value = 7289.8335
$ws.Range("J51").Value = 7907.091
$ws.Range("L51").Value = 7907.091
$ws.Range("N51").Value = -9379.091

$ws.Range("H58").Value = 2244.8845
$ws.Range("I58").Value = 2131.1333
$ws.Range("J58").Value = 2400
$ws.Range("K58").Value = 2131.1333
$ws.Range("L58").Value = 2400
$ws.Range("M58").Value = -1928.1333
$ws.Range("N58").Value = -2806

$ws.Range("H60").Value = 15020.75
$ws.Range("J60").Value = 15020.75
$ws.Range("L60").Value = 15020.75
$ws.Range("N60").Value = -16042.75

$ws.Range("H61").Value = 7289.8335
$ws.Range("J61").Value = 7907.091
$ws.Range("L61").Value = 7907.091
$ws.Range("N61").Value = -8603.091

$ws.Range("H68").Value = 17277.05
$ws.Range("J68").Value = 17277.05
$ws.Range("L68").Value = 17277.05
$ws.Range("N68").Value = -18775.05

$ws.Range("H71").Value = 17277.05
$ws.Range("J71").Value = 17277.05
$ws.Range("L71").Value = 51831.14999999999
$ws.Range("N71").Value = -59319.14999999999

$ws.Range("H74").Value = 40900
$ws.Range("J74").Value = 40900
$ws.Range("L74").Value = 40900
$ws.Range("N74").Value = -42648

$ws.Range("H77").Value = 40900
$ws.Range("J77").Value = 40900
$ws.Range("L77").Value = 122700
$ws.Range("N77").Value = -131436

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = $null

$ws.Range("H132").Value = 3204.3076
$ws.Range("I132").Value = 3184.8823
$ws.Range("J132").Value = 3241
$ws.Range("K132").Value = 9554.6469
$ws.Range("L132").Value = 9723
$ws.Range("M132").Value = -7024.6469
$ws.Range("N132").Value = -14783

$ws.Range("H136").Value = 2244.8845
$ws.Range("I136").Value = 2131.1333
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 6393.3999
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -3843.3999
$ws.Range("N136").Value = -12300

$ws.Range("H141").Value = 104793.14
$ws.Range("J141").Value = 92666.664
$ws.Range("L141").Value = 92666.664
$ws.Range("N141").Value = -103026.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 142.25
$ws.Range("J38").Value = 161.6
$ws.Range("L38").Value = 484.8
$ws.Range("N38").Value = -1178.8

$ws.Range("H107").Value = 1378.1666
$ws.Range("I107").Value = 801.8570999999999
$ws.Range("J107").Value = 1744.909
$ws.Range("K107").Value = 2405.5713
$ws.Range("L107").Value = 5234.727000000001
$ws.Range("M107").Value = -485.5712999999996
$ws.Range("N107").Value = -9074.727000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3111.4285
$ws.Range("I126").Value = 2896
$ws.Range("J126").Value = 3650
$ws.Range("K126").Value = 8688
$ws.Range("L126").Value = 10950
$ws.Range("M126").Value = -6218
$ws.Range("N126").Value = -15890

$ws.Range("H137").Value = 68000
$ws.Range("J137").Value = 68000
$ws.Range("L137").Value = 68000
$ws.Range("N137").Value = -78200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 73859.21000000001
$ws.Range("I40").Value = 168779.83
$ws.Range("J40").Value = 2668.75
$ws.Range("K40").Value = 168779.83
$ws.Range("L40").Value = 2668.75
$ws.Range("M40").Value = -168643.83
$ws.Range("N40").Value = -2940.75

$ws.Range("H132").Value = 4124.793
$ws.Range("I132").Value = 5280.625
$ws.Range("K132").Value = 15841.875
$ws.Range("M132").Value = -13311.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 16230
$ws.Range("J123").Value = 16230
$ws.Range("L123").Value = 16230
$ws.Range("N123").Value = -26030

$ws.Range("H126").Value = 1476.6666
$ws.Range("I126").Value = 1960
$ws.Range("J126").Value = 993.3333
$ws.Range("K126").Value = 5880
$ws.Range("L126").Value = 2979.9999
$ws.Range("M126").Value = -3410
$ws.Range("N126").Value = -7919.9999

$ws.Range("H132").Value = 23494.541
$ws.Range("I132").Value = 3152.2
$ws.Range("J132").Value = 57398.445
$ws.Range("K132").Value = 9456.599999999999
$ws.Range("L132").Value = 172195.335
$ws.Range("M132").Value = -6926.599999999999
$ws.Range("N132").Value = -177255.335

$ws.Range("H136").Value = 3969.0852
$ws.Range("I136").Value = 6822
$ws.Range("J136").Value = 2352.4333
$ws.Range("K136").Value = 20466
$ws.Range("L136").Value = 7057.2999
$ws.Range("M136").Value = -17916
$ws.Range("N136").Value = -12157.2999
